$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.444.08"
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("D3").Value = "'1.852.83"
$ws.Range("E3").Value = "  +1.15%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'233.47"
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("D7").Value = "'0.4751"
$ws.Range("E7").Value = "  +2.24%  "
$ws.Range("D8").Value = "'0.2754"
$ws.Range("E8").Value = "  +2.28%  "
$ws.Range("D9").Value = "'0.06342"
$ws.Range("E9").Value = "  +1.35%  "
$ws.Range("D10").Value = "'1.970.76"
$ws.Range("E10").Value = "  +7.52%  "
$ws.Range("D11").Value = "'17.83"
$ws.Range("E11").Value = "  +11.18%  "
$ws.Range("D12").Value = "'0.07461"
$ws.Range("E12").Value = "  +1.07%  "
$ws.Range("D13").Value = "'4.960"
$ws.Range("E13").Value = "  +1.42%  "
$ws.Range("D14").Value = "'84.76"
$ws.Range("E14").Value = "  +1.84%  "
$ws.Range("D15").Value = "'0.6249"
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("D16").Value = "'30.403.58"
$ws.Range("E16").Value = "  +1.05%  "
$ws.Range("D17").Value = "'245.94"
$ws.Range("E17").Value = "  +8.33%  "
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").Value = "'12.67"
$ws.Range("E19").Value = "  +2.35%  "
$ws.Range("D20").Value = "'0.000007334"
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "'4.919"
$ws.Range("E22").Value = "  +1.67%  "
$ws.Range("D23").Value = "'5.917"
$ws.Range("E23").Value = "  +0.89%  "
$ws.Range("D24").Value = "'164.40"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").Value = "'9.064"
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("D26").Value = "'17.98"
$ws.Range("E26").Value = "  +1.89%  "
$ws.Range("D27").Value = "'1.876"
$ws.Range("E27").Value = "  +1.87%  "
$ws.Range("D28").Value = "'0.1028"
$ws.Range("E28").Value = "  +2.07%  "
$ws.Range("E29").Value = "  -1.40%  "
$ws.Range("D30").Value = "'4.040"
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("D31").Value = "'3.831"
$ws.Range("E31").Value = "  +1.85%  "
$ws.Range("D32").Value = "'0.04833"
$ws.Range("E32").Value = "  +1.19%  "
$ws.Range("D33").Value = "'1.129"
$ws.Range("E33").Value = "  +0.46%  "
$ws.Range("D34").Value = "'0.6987"
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("D35").Value = "'2.705"
$ws.Range("E35").Value = "  +0.66%  "
$ws.Range("D36").Value = "'0.01905"
$ws.Range("E36").Value = "  +5.56%  "
$ws.Range("D37").Value = "'2.681"
$ws.Range("E37").Value = "  +2.67%  "
$ws.Range("D38").Value = "'0.8778"
$ws.Range("E38").Value = "  -1.89%  "
$ws.Range("D39").Value = "'1.990"
$ws.Range("E39").Value = "  +3.55%  "
$ws.Range("D40").Value = "'106.67"
$ws.Range("E40").Value = "  +3.78%  "
$ws.Range("D41").Value = "'0.9999"
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").Value = "'0.4056"
$ws.Range("E42").Value = "  +1.64%  "
$ws.Range("D43").Value = "'5.499"
$ws.Range("E43").Value = "  +0.20%  "
$ws.Range("D44").Value = "'7.170"
$ws.Range("E44").Value = "  +3.26%  "
$ws.Range("D45").Value = "'63.32"
$ws.Range("E45").Value = "  +6.45%  "
$ws.Range("D46").Value = "'0.1200"
$ws.Range("E46").Value = "  +0.91%  "
$ws.Range("D47").Value = "'34.00"
$ws.Range("D48").Value = "'8.533"
$ws.Range("E48").Value = "  +1.38%  "
$ws.Range("D49").Value = "'0.05499"
$ws.Range("E49").Value = "  -0.43%  "
$ws.Range("E50").Value = "  -0.88%  "
$ws.Range("D51").Value = "'0.3686"
$ws.Range("E51").Value = "  +1.96%  "
